$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete trailing rows 27-29 (content removed by the edit; shifts nothing else since
# rows 1-26 are rewritten explicitly below)
$ws.Rows.Item(27).Delete()
$ws.Rows.Item(27).Delete()
$ws.Rows.Item(27).Delete()

# Row 10
$ws.Range("B10").Value = '144651 - Antonio Fernando Sartori'
$ws.Range("C10").Value = '144651 - Antonio Fernando Sartori'

# Row 13
$ws.Range("A13").Value = 'Programa resumido:'
$ws.Rows.Item(13).RowHeight = 60

# Row 14
$ws.Range("B14").ClearContents()
$ws.Range("C14").ClearContents()
$ws.Range("A14").Value = 'Short syllabus:'
$ws.Rows.Item(14).RowHeight = 60

# Row 15
$ws.Range("A15").Value = 'Programa:'
$ws.Range("B15").Value = '5840730 - Antonio Jefferson da Silva Machado'
$ws.Range("C15").Value = '5840730 - Antonio Jefferson da Silva Machado'
$ws.Rows.Item(15).RowHeight = 120

# Row 16
$ws.Range("B16").ClearContents()
$ws.Range("C16").ClearContents()
$ws.Range("A16").Value = 'Syllabus:'
$ws.Rows.Item(16).RowHeight = 120

# Row 17
$ws.Range("A17").Value = 'Avaliação:'
$ws.Rows.Item(17).AutoFit()

# Row 18
$ws.Range("A18").Value = 'Método:'
$ws.Range("B18").Value = '519033 - Carlos Yujiro Shigue'
$ws.Range("C18").Value = '519033 - Carlos Yujiro Shigue'
$ws.Rows.Item(18).RowHeight = 60

# Row 19
$ws.Range("A19").Value = 'Critério:'
$ws.Range("B19").Value = 'Aulas expositivas, práticas, seminários e exercícios.'
$ws.Range("C19").Value = 'Aulas expositivas, práticas, seminários e exercícios.'
$ws.Rows.Item(19).RowHeight = 60

# Row 20
$ws.Range("A20").Value = 'Norma de recuperação:'
$ws.Range("B20").Value = 'Média das notas de provas, relatórios e apresentações.'
$ws.Range("C20").Value = 'Média das notas de provas, relatórios e apresentações.'
$ws.Rows.Item(20).RowHeight = 60

# Row 21
$ws.Range("A21").Value = 'Bibliografia:'
$ws.Range("B21").Value = 'Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação'
$ws.Range("C21").Value = 'Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação'
$ws.Rows.Item(21).RowHeight = 120

# Row 22
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()
$ws.Range("A22").Value = 'Requisitos:'
$ws.Rows.Item(22).AutoFit()

# Row 23
$ws.Range("A23").ClearContents()
$ws.Range("B23").Value = "LOM3206 -  Eletrônica  (Requisito)`n"
$ws.Range("C23").Value = "LOM3206 -  Eletrônica  (Requisito)`n"
$ws.Rows.Item(23).RowHeight = 30

# Row 24
$ws.Range("A24").ClearContents()
$ws.Range("B24").Value = "LOM3215 -  Física do Estado Sólido  (Requisito)`n"
$ws.Range("C24").Value = "LOM3215 -  Física do Estado Sólido  (Requisito)`n"
$ws.Rows.Item(24).RowHeight = 30

# Row 25
$ws.Range("A25").ClearContents()
$ws.Range("B25").Value = "LOM3231 -  Métodos Experimentais da Física IV  (Indicação de Conjunto)`n"
$ws.Range("C25").Value = "LOM3231 -  Métodos Experimentais da Física IV  (Indicação de Conjunto)`n"
$ws.Rows.Item(25).RowHeight = 30

# Row 26
$ws.Range("B26").Value = "LOM3234 -  Óptica Física  (Requisito)`n"
$ws.Range("C26").Value = "LOM3234 -  Óptica Física  (Requisito)`n"
